$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old column H ("Total Ulhi"), shifting
# everything from H onward one column to the right (H->I, I->J, ... O->P, P->Q)
$ws.Columns.Item(8).Insert()
$ws.Range("H1").Value = "Underlying Health Issues"

# --- Update existing row 2 (was John Doe, now Chad Oliver) ---
$ws.Range("B2").Value = "Chad"
$ws.Range("C2").Value = "Oliver"
$ws.Range("D2").Value = "chadoliver017@gmail.com"
$ws.Range("E2").Value = 96.8
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = "Pains"
$ws.Range("H2").Value = "Asthma"
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = $false
$ws.Range("P2").Value = $false
$ws.Range("Q2").Value = "Not at Risk"

# --- Row 3 (new) ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Chad"
$ws.Range("C3").Value = "Oliver"
$ws.Range("D3").Value = "chadoliver017@gmail.com"
$ws.Range("E3").Value = 100.4
$ws.Range("F3").Value = 20
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "Asthma"
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = $false
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = "Not at Risk"

# --- Row 4 (new) ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Chad"
$ws.Range("C4").Value = "Oliver"
$ws.Range("D4").Value = "chadoliver017@gmail.com"
$ws.Range("E4").Value = 100.4
$ws.Range("F4").Value = 20
$ws.Range("G4").Value = "blank"
$ws.Range("H4").Value = "Asthma"
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = $false
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = "Not at Risk"

# --- Row 5 (new) ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Chad"
$ws.Range("C5").Value = "Oliver"
$ws.Range("D5").Value = "chadoliver017@gmail.com"
$ws.Range("E5").Value = 109.4
$ws.Range("F5").Value = 21
$ws.Range("G5").Value = "Chest Pain,Loss of Speech"
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = $false
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = "Very High Risk"
